$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.222.66"
$ws.Range("E2").Value = "  +3.15%  "
$ws.Range("D3").Value = "1.898.73"
$ws.Range("E3").Value = "  -0.02%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.66"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.31%  "
$ws.Range("E6").Value = "  -0.18%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5174"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.46%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4023"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.33%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08452"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.75"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.60%  "
$ws.Range("E11").Value = "  +0.09%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "23.57"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +13.47%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.438"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.97%  "
$ws.Range("D14").Value = "1.912.07"
$ws.Range("E14").Value = "  +0.78%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.338"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.17%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.003"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.10%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "94.85"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.74%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001113"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.60%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.28"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.37%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.002"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.966"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.02%  "
$ws.Range("D23").Value = "30.245.92"
$ws.Range("E23").Value = "  +3.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.27"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.86%  "
$ws.Range("E25").Value = "  +0.20%  "
$ws.Range("D26").Value = "2.122.93"
$ws.Range("E26").Value = "  +0.47%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.74"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.84%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "161.28"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.33%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.398"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.61%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "129.53"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.92%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.090"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.69%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1058"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.04%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.011"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.68%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.740"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.27%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02495"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.88%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06569"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.11%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2206"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.52%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.238"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.28%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.219"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.25%  "
$ws.Range("E40").Value = "  +4.95%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.775"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.23%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6495"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.06%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.234"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.24%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6101"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.51%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.25"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.49%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.711"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.98%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.055"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.29%  "
$ws.Range("E48").Value = "  +0.89%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "124.77"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.10%  "
$ws.Range("E50").Value = "  -0.96%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "79.18"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.81%  "
